$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Regenerated handoff report: the localized file's GUID-based name changed
# from 9c541159-abff-48bf-9958-f3030b7106f1 to b742e51e-0df9-44be-a16e-1a022713b4da,
# and the handoff/handback timestamps moved forward a few seconds.
# ---------------------------------------------------------------------------

$oldGuid = "9c541159-abff-48bf-9958-f3030b7106f1"
$newGuid = "b742e51e-0df9-44be-a16e-1a022713b4da"

$oldHash = "5f636792f8d2d3fd32bb7bea717ce529d18376a5"
$newHash = "25b5561527b3d21c1c8e4884f0d0b37954a58ac1"

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cc780461921522db2dcab7fe6fbc0dcfdbc373d6/e2e/$oldGuid.md"

# --- Overview sheet --------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkTarget, "", "", "e2e\$newGuid.md") | Out-Null

$wsOverview.Range("G2").Value = "2016-08-15 20:53:54"

# --- zh-cn sheet -------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkTarget, "", "", "$newGuid.md") | Out-Null

$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-15 20:53:50"

# --- de-de sheet -------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkTarget, "", "", "$newGuid.md") | Out-Null

$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-15 20:53:54"
